# Rename the PARMA-related linked fields to PAMA across the workbook
# (survey, queries, model sheets) and update the related selection / view
# state, per commit "Fixed linked tables icluding new hosp. and fal."

$wb = $excel.ActiveWorkbook

# --- survey sheet -----------------------------------------------------
$survey = $wb.Worksheets.Item("survey")
$survey.Range("E20").Value = "linked_pama"
$survey.Range("F20").Value = "PAMA"
$survey.Range("E21").Value = "linked_dapama"
$survey.Range("F21").Value = "DAPAMA"
$survey.Range("E22").Value = "linked_idpama"
$survey.Range("F22").Value = "IDPAMA"

# --- queries sheet ------------------------------------------------------
$queries = $wb.Worksheets.Item("queries")
$queries.Range("A11").Value = "linked_pama"
$queries.Range("I11").Value = "PAMA"
$queries.Range("A12").Value = "linked_dapama"
$queries.Range("I12").Value = "DAPAMA"
$queries.Range("A13").Value = "linked_idpama"
$queries.Range("I13").Value = "IDPAMA"
$queries.Range("G2").Value = "{REGIDC: data('REGIDC'), comsup: data('COMSUP'), pama: data('PAMA'), moma: data('MOMA'), REGID: data('REGID')}"

# --- model sheet ----------------------------------------------------
$model = $wb.Worksheets.Item("model")
$model.Range("A23").Value = "PAMA"
$model.Range("A24").Value = "DAPAMA"
$model.Range("A25").Value = "IDPAMA"

# --- selection / view state -------------------------------------------
$survey.Activate()
$survey.Range("F20").Select()

$model.Activate()
$model.Application.ActiveWindow.ScrollRow = 2
$model.Range("A25").Select()
